# Update the "dSF" column (F) values to reflect re-pulled / re-pushed data
# and recomputed mean, per commit message "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    6  = 6
    7  = -3
    8  = 6
    9  = 2
    10 = -2
    11 = -2
    12 = -1
    13 = -1
    14 = 1
    15 = -3
    16 = -3
    17 = 3
    18 = 3
    19 = 11
    20 = 1
    21 = -7
    22 = -2
    23 = 2
    24 = -2
    25 = -1
    26 = -1
    27 = 5
    28 = -2
    29 = 3
    30 = 7
    31 = -4
    32 = 2
    33 = 5
    34 = -2
    35 = 4
    37 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
